{"js": "// Replace the 25 two-digit multiplication problems in the practice table\n// with the values from the authored revision. Each old problem string is\n// unique in the document, so an exact, case-sensitive search-and-replace\n// on the document body is sufficient and unambiguous.\nconst replacements = [\n  [\"37\u00d749=\", \"27\u00d779=\"],\n  [\"69\u00d749=\", \"77\u00d753=\"],\n  [\"49\u00d782=\", \"61\u00d793=\"],\n  [\"45\u00d721=\", \"62\u00d766=\"],\n  [\"95\u00d736=\", \"23\u00d765=\"],\n  [\"99\u00d740=\", \"26\u00d733=\"],\n  [\"40\u00d795=\", \"63\u00d760=\"],\n  [\"79\u00d754=\", \"15\u00d741=\"],\n  [\"98\u00d771=\", \"71\u00d718=\"],\n  [\"50\u00d785=\", \"13\u00d774=\"],\n  [\"28\u00d725=\", \"69\u00d783=\"],\n  [\"37\u00d728=\", \"92\u00d734=\"],\n  [\"34\u00d738=\", \"89\u00d741=\"],\n  [\"75\u00d779=\", \"98\u00d748=\"],\n  [\"96\u00d742=\", \"42\u00d767=\"],\n  [\"60\u00d713=\", \"12\u00d760=\"],\n  [\"43\u00d744=\", \"83\u00d797=\"],\n  [\"53\u00d778=\", \"87\u00d747=\"],\n  [\"46\u00d771=\", \"28\u00d780=\"],\n  [\"92\u00d767=\", \"77\u00d722=\"],\n  [\"86\u00d724=\", \"76\u00d766=\"],\n  [\"27\u00d730=\", \"24\u00d730=\"],\n  [\"45\u00d782=\", \"72\u00d741=\"],\n  [\"59\u00d766=\", \"67\u00d720=\"],\n  [\"81\u00d776=\", \"13\u00d761=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit multiplication problems in the practice table\n# with the values from the authored revision. Each old problem string is\n# unique in the document, so a simple Find/Replace-All per pair is\n# unambiguous and safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"37\u00d749=\", \"27\u00d779=\"),\n    @(\"69\u00d749=\", \"77\u00d753=\"),\n    @(\"49\u00d782=\", \"61\u00d793=\"),\n    @(\"45\u00d721=\", \"62\u00d766=\"),\n    @(\"95\u00d736=\", \"23\u00d765=\"),\n    @(\"99\u00d740=\", \"26\u00d733=\"),\n    @(\"40\u00d795=\", \"63\u00d760=\"),\n    @(\"79\u00d754=\", \"15\u00d741=\"),\n    @(\"98\u00d771=\", \"71\u00d718=\"),\n    @(\"50\u00d785=\", \"13\u00d774=\"),\n    @(\"28\u00d725=\", \"69\u00d783=\"),\n    @(\"37\u00d728=\", \"92\u00d734=\"),\n    @(\"34\u00d738=\", \"89\u00d741=\"),\n    @(\"75\u00d779=\", \"98\u00d748=\"),\n    @(\"96\u00d742=\", \"42\u00d767=\"),\n    @(\"60\u00d713=\", \"12\u00d760=\"),\n    @(\"43\u00d744=\", \"83\u00d797=\"),\n    @(\"53\u00d778=\", \"87\u00d747=\"),\n    @(\"46\u00d771=\", \"28\u00d780=\"),\n    @(\"92\u00d767=\", \"77\u00d722=\"),\n    @(\"86\u00d724=\", \"76\u00d766=\"),\n    @(\"27\u00d730=\", \"24\u00d730=\"),\n    @(\"45\u00d782=\", \"72\u00d741=\"),\n    @(\"59\u00d766=\", \"67\u00d720=\"),\n    @(\"81\u00d776=\", \"13\u00d761=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
